# Auto-generated Excel COM-interop script
# Applies numeric value updates (from the authoritative commit diff)
# to the profit-tracking columns (H-N) of several worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")

# Row 30
$ws.Range("H30").Value = 2256.6667
$ws.Range("I30").Value = 9
$ws.Range("J30").Value = 3380.5
$ws.Range("K30").Value = 9
$ws.Range("L30").Value = 3380.5
$ws.Range("M30").Value = 141
$ws.Range("N30").Value = -3680.5

# Row 63
$ws.Range("H63").Value = 2498
$ws.Range("I63").Value = 1297.8
$ws.Range("K63").Value = 1297.8
$ws.Range("M63").Value = -611.8

# Row 66
$ws.Range("H66").Value = 2498
$ws.Range("I66").Value = 1297.8
$ws.Range("K66").Value = 6489
$ws.Range("M66").Value = -3057

# Row 135
$ws.Range("H135").Value = 59500
$ws.Range("J135").Value = 59500
$ws.Range("L135").Value = 59500
$ws.Range("N135").Value = -69640


$ws = $wb.Worksheets.Item("BSM")

# Row 8
$ws.Range("H8").Value = 109
$ws.Range("I8").Value = 109
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 109
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 31
$ws.Range("N8").Value = ""

# Row 10
$ws.Range("H10").Value = 998.3333
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 998.3333
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 998.3333
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = -1278.3333

# Row 11
$ws.Range("H11").Value = 710
$ws.Range("J11").Value = 2000
$ws.Range("L11").Value = 2000
$ws.Range("N11").Value = -2280

# Row 12
$ws.Range("H12").Value = 670
$ws.Range("I12").Value = 50
$ws.Range("J12").Value = 1083.3334
$ws.Range("K12").Value = 50
$ws.Range("L12").Value = 1083.3334
$ws.Range("M12").Value = 118
$ws.Range("N12").Value = -1419.3334

# Row 23
$ws.Range("H23").Value = 3708.5
$ws.Range("J23").Value = 3708.5
$ws.Range("L23").Value = 3708.5
$ws.Range("N23").Value = -4274.5

# Row 81
$ws.Range("H81").Value = 14112
$ws.Range("J81").Value = 14112
$ws.Range("L81").Value = 14112
$ws.Range("N81").Value = -16234

# Row 84
$ws.Range("H84").Value = 14112
$ws.Range("J84").Value = 14112
$ws.Range("L84").Value = 42336
$ws.Range("N84").Value = -52944

# Row 97
$ws.Range("H97").Value = 24400
$ws.Range("I97").Value = 24400
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 24400
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -23409
$ws.Range("N97").Value = ""

# Row 98
$ws.Range("H98").Value = 10000
$ws.Range("I98").Value = 10000
$ws.Range("K98").Value = 10000
$ws.Range("M98").Value = -7005

# Row 99
$ws.Range("H99").Value = 2330.9
$ws.Range("I99").Value = 2309.8333
$ws.Range("K99").Value = 2309.8333
$ws.Range("M99").Value = -811.8332999999998

# Row 102
$ws.Range("H102").Value = 17811.25
$ws.Range("I102").Value = 17811.25
$ws.Range("K102").Value = 17811.25
$ws.Range("M102").Value = -14566.25


$ws = $wb.Worksheets.Item("CRP")

# Row 3
$ws.Range("H3").Value = 733
$ws.Range("I3").Value = 733
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 733
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -620
$ws.Range("N3").Value = ""

# Row 52
$ws.Range("H52").Value = 49186.332
$ws.Range("I52").Value = 50000
$ws.Range("J52").Value = 48779.5
$ws.Range("K52").Value = 50000
$ws.Range("L52").Value = 48779.5
$ws.Range("M52").Value = -49706
$ws.Range("N52").Value = -49367.5

# Row 132
$ws.Range("H132").Value = 7719.2188
$ws.Range("I132").Value = 6039.077
$ws.Range("K132").Value = 18117.231
$ws.Range("M132").Value = -15587.231

# Row 134
$ws.Range("H134").Value = 1905.2106
$ws.Range("I134").Value = 1696.8125
$ws.Range("K134").Value = 5090.4375
$ws.Range("M134").Value = -2555.4375

# Row 141
$ws.Range("H141").Value = 383333
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 383333
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 383333
$ws.Range("M141").Value = ""
$ws.Range("N141").Value = -393693


$ws = $wb.Worksheets.Item("CUL")

# Row 97
$ws.Range("H97").Value = 623.75
$ws.Range("I97").Value = 750
$ws.Range("J97").Value = 497.5
$ws.Range("K97").Value = 2250
$ws.Range("L97").Value = 1492.5
$ws.Range("M97").Value = -1754
$ws.Range("N97").Value = -2484.5

# Row 121
$ws.Range("H121").Value = 1066.6666
$ws.Range("I121").Value = 100
$ws.Range("K121").Value = 300
$ws.Range("M121").Value = 1010


$ws = $wb.Worksheets.Item("GSM")

# Row 12
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").Value = ""

# Row 13
$ws.Range("H13").Value = 217
$ws.Range("I13").Value = 100.5
$ws.Range("J13").Value = 450
$ws.Range("K13").Value = 100.5
$ws.Range("L13").Value = 450
$ws.Range("M13").Value = 38.5
$ws.Range("N13").Value = -728

# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = ""

# Row 22
$ws.Range("H22").Value = 1710.3334
$ws.Range("I22").Value = 1710.3334
$ws.Range("K22").Value = 1710.3334
$ws.Range("M22").Value = -1181.3334

# Row 25
$ws.Range("H25").Value = 19999
$ws.Range("J25").Value = 19999
$ws.Range("L25").Value = 19999
$ws.Range("N25").Value = -21057

# Row 27
$ws.Range("H27").Value = 50999.5
$ws.Range("I27").Value = 100000
$ws.Range("J27").Value = 1999
$ws.Range("K27").Value = 100000
$ws.Range("L27").Value = 1999
$ws.Range("M27").Value = -99834
$ws.Range("N27").Value = -2331

# Row 30
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = ""

# Row 113
$ws.Range("H113").Value = 1811
$ws.Range("I113").Value = 1811
$ws.Range("K113").Value = 1811
$ws.Range("M113").Value = 359


$ws = $wb.Worksheets.Item("LTW")

# Row 64
$ws.Range("H64").Value = 12500
$ws.Range("J64").Value = 12500
$ws.Range("L64").Value = 12500
$ws.Range("N64").Value = -12950

# Row 67
$ws.Range("H67").Value = 12500
$ws.Range("J67").Value = 12500
$ws.Range("L67").Value = 12500
$ws.Range("N67").Value = -14060

# Row 82
$ws.Range("H82").Value = 3822.2222
$ws.Range("J82").Value = 4150
$ws.Range("L82").Value = 4150
$ws.Range("N82").Value = -4872

# Row 85
$ws.Range("H85").Value = 3822.2222
$ws.Range("J85").Value = 4150
$ws.Range("L85").Value = 4150
$ws.Range("N85").Value = -6646

# Row 122
$ws.Range("H122").Value = 6024.7334
$ws.Range("I122").Value = 5438.2
$ws.Range("K122").Value = 16314.6
$ws.Range("M122").Value = -13864.6


$ws = $wb.Worksheets.Item("WVR")

# Row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = ""

# Row 23
$ws.Range("H23").Value = 347.66666
$ws.Range("I23").Value = 300
$ws.Range("J23").Value = 443
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 443
$ws.Range("M23").Value = -71
$ws.Range("N23").Value = -901

# Row 63
$ws.Range("H63").Value = 19999
$ws.Range("J63").Value = 19999
$ws.Range("L63").Value = 19999
$ws.Range("N63").Value = -21247

# Row 66
$ws.Range("H66").Value = 19999
$ws.Range("J66").Value = 19999
$ws.Range("L66").Value = 59997
$ws.Range("N66").Value = -66237

# Row 107
$ws.Range("H107").Value = 522.75
$ws.Range("I107").Value = 522.75
$ws.Range("K107").Value = 1568.25
$ws.Range("M107").Value = 351.75

# Row 122
$ws.Range("H122").Value = 1143.2222
$ws.Range("J122").Value = 649.5
$ws.Range("L122").Value = 1948.5
$ws.Range("N122").Value = -6848.5

